# === Workbook / sheet setup ===
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add "Answers" sheet right after HelpBoxText
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Answers"

# Add "QuestionsTitles" sheet right after Answers
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "QuestionsTitles"

# === Populate "Answers" sheet (column A) ===
$ws2.Range("A1").Value = "BMW Financial Services offers a variety of easy, convenient payment methods.`nOnline with your My BMW account`nEasily manage your account online with My BMW. From here, you can make one-time electronic payments, enroll in EasyPay automatic payments, sign up for Paperless Statements, view your Account Statements, and more.`nLog in or create your account today at mybmw.bmwusa.com`nPay By Phone`nCall 800-578-5000 and make a one-time payment with your checking/savings account or debit card for the same day, or schedule your payment for a future date.`nPay By Mail`nTear off and return the bottom portion of your Account Statement with your payment. (If you're signed up for paperless statements, print the return portion). Make checks payable to BMW Financial Services and include your 10-digit BMW Financial Services account number."
$ws2.Range("A2").Value = "Account Statements are delivered approximately 14 days before a payment is due. You'll receive a notification by mail, or by email if you're enrolled in paperless statements.`nYou can also view your Account Statement any time by signing in to My BMW."
$ws2.Range("A3").Value = "Payments are first credited to unpaid finance charges, then to the outstanding principal balance, then to any outstanding fees. Your Account Statement will break down how your payments have been allocated.`nAmortization schedules are always available through your My BMW account."
$ws2.Range("A4").Value = "While debit cards are acceptable, we unfortunately cannot accept credit cards for regular monthly payments."
$ws2.Range("A5").Value = "Grace periods – the amount of time between a due date and the assessment of late fees – are regulated by state governments. Please refer to the Account Summary section of your Account Statement to determine if you are eligible for a grace period."
$ws2.Range("A6").Value = "Yes. If you pay more than the Total Amount Due, the extra payment will be applied to your principal balance. This may reduce the interest you pay over the life of your account and may reduce your final payment or shorten the term of your financing agreement."
$ws2.Range("A7").Value = "Payment credits are applied on your scheduled due date and will be reflected on your Account Statement."
$ws2.Range("A8").Value = "Late payments (more than 29 days past due), missed payments, or other defaults on your account may be reflected on your credit report. In accordance with Federal law, you are hereby notified that a negative credit report reflecting on your credit records may be submitted to a credit reporting agency if you fail to fulfill the terms of your credit obligation."

# Wrap text for all used cells on Answers
$ws2.Range("A1:A8").WrapText = $true

# Row heights on Answers
$ws2.Rows.Item(1).RowHeight = 152.25
$ws2.Rows.Item(2).RowHeight = 45
$ws2.Rows.Item(3).RowHeight = 45
$ws2.Rows.Item(5).RowHeight = 30
$ws2.Rows.Item(6).RowHeight = 30
$ws2.Rows.Item(8).RowHeight = 45

# Column width on Answers
$ws2.Columns.Item(1).ColumnWidth = 125.28515625

# === Populate "QuestionsTitles" sheet (column A) ===
$ws3.Range("A1").Value = "1`nWhat payment options are available?"
$ws3.Range("A2").Value = "2`nWhen will I receive my monthly Account Statements?"
$ws3.Range("A3").Value = "3`nHow are my payments applied?"
$ws3.Range("A4").Value = "4`nCan I make a payment with my debit or credit card?"
$ws3.Range("A5").Value = "5`nIs there a grace period for late charges?"
$ws3.Range("A6").Value = "6`nCan I make a payment greater than my regular scheduled monthly payment?"
$ws3.Range("A7").Value = "7`nHow are payment credits applied?"
$ws3.Range("A8").Value = "8`nWhat happens if I miss a payment?"

# Wrap text for all used cells on QuestionsTitles
$ws3.Range("A1:A8").WrapText = $true

# Row heights on QuestionsTitles
$ws3.Rows.Item(1).RowHeight = 39
$ws3.Rows.Item(2).RowHeight = 45
$ws3.Rows.Item(3).RowHeight = 30
$ws3.Rows.Item(4).RowHeight = 45
$ws3.Rows.Item(5).RowHeight = 30
$ws3.Rows.Item(6).RowHeight = 45
$ws3.Rows.Item(7).RowHeight = 30
$ws3.Rows.Item(8).RowHeight = 30

# Column width on QuestionsTitles
$ws3.Columns.Item(1).ColumnWidth = 46.28515625

# === Selection / active-sheet / view state ===
# HelpBoxText keeps its existing A3 selection; it simply stops being the active
# sheet (tabSelected) once another sheet is activated below.

# Answers: default selection (A1)
$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null

# QuestionsTitles ends up the active / tab-selected sheet, with A1:A8 selected
$ws3.Activate() | Out-Null
$ws3.Range("A1:A8").Select() | Out-Null
